$d = $word.ActiveDocument

$d.Content.Find.Execute("92-74=18", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=94", 2) | Out-Null
$d.Content.Find.Execute("24+65=89", $true, $false, $false, $false, $false, $true, 1, $false, "82-73=9", 2) | Out-Null
$d.Content.Find.Execute("26+43=69", $true, $false, $false, $false, $false, $true, 1, $false, "85-62=23", 2) | Out-Null
$d.Content.Find.Execute("53+4=57", $true, $false, $false, $false, $false, $true, 1, $false, "52-12=40", 2) | Out-Null
$d.Content.Find.Execute("26-11=15", $true, $false, $false, $false, $false, $true, 1, $false, "99-59=40", 2) | Out-Null
$d.Content.Find.Execute("42-6=36", $true, $false, $false, $false, $false, $true, 1, $false, "53+1=54", 2) | Out-Null
$d.Content.Find.Execute("26+11=37", $true, $false, $false, $false, $false, $true, 1, $false, "86-78=8", 2) | Out-Null
$d.Content.Find.Execute("51-2=49", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=91", 2) | Out-Null
$d.Content.Find.Execute("55+23=78", $true, $false, $false, $false, $false, $true, 1, $false, "72-2=70", 2) | Out-Null
$d.Content.Find.Execute("72-14=58", $true, $false, $false, $false, $false, $true, 1, $false, "92-9=83", 2) | Out-Null
$d.Content.Find.Execute("14+80=94", $true, $false, $false, $false, $false, $true, 1, $false, "86-80=6", 2) | Out-Null
$d.Content.Find.Execute("48+28=76", $true, $false, $false, $false, $false, $true, 1, $false, "74-48=26", 2) | Out-Null
$d.Content.Find.Execute("51+5=56", $true, $false, $false, $false, $false, $true, 1, $false, "9+34=43", 2) | Out-Null
$d.Content.Find.Execute("32-2=30", $true, $false, $false, $false, $false, $true, 1, $false, "8-1=7", 2) | Out-Null
$d.Content.Find.Execute("32+50=82", $true, $false, $false, $false, $false, $true, 1, $false, "14+57=71", 2) | Out-Null
$d.Content.Find.Execute("15+10=25", $true, $false, $false, $false, $false, $true, 1, $false, "61-39=22", 2) | Out-Null
$d.Content.Find.Execute("74-21=53", $true, $false, $false, $false, $false, $true, 1, $false, "3+76=79", 2) | Out-Null
$d.Content.Find.Execute("90-11=79", $true, $false, $false, $false, $false, $true, 1, $false, "20+17=37", 2) | Out-Null
$d.Content.Find.Execute("79+4=83", $true, $false, $false, $false, $false, $true, 1, $false, "64-42=22", 2) | Out-Null
$d.Content.Find.Execute("6+38=44", $true, $false, $false, $false, $false, $true, 1, $false, "81-23=58", 2) | Out-Null
$d.Content.Find.Execute("13-10=3", $true, $false, $false, $false, $false, $true, 1, $false, "30+23=53", 2) | Out-Null
$d.Content.Find.Execute("62+31=93", $true, $false, $false, $false, $false, $true, 1, $false, "1+75=76", 2) | Out-Null
$d.Content.Find.Execute("24-4=20", $true, $false, $false, $false, $false, $true, 1, $false, "84-29=55", 2) | Out-Null
$d.Content.Find.Execute("94-7=87", $true, $false, $false, $false, $false, $true, 1, $false, "93-10=83", 2) | Out-Null
$d.Content.Find.Execute("56+17=73", $true, $false, $false, $false, $false, $true, 1, $false, "89-20=69", 2) | Out-Null
$d.Content.Find.Execute("41+47=88", $true, $false, $false, $false, $false, $true, 1, $false, "57-31=26", 2) | Out-Null
$d.Content.Find.Execute("50+20=70", $true, $false, $false, $false, $false, $true, 1, $false, "20+52=72", 2) | Out-Null
$d.Content.Find.Execute("49-28=21", $true, $false, $false, $false, $false, $true, 1, $false, "20+46=66", 2) | Out-Null
$d.Content.Find.Execute("85+5=90", $true, $false, $false, $false, $false, $true, 1, $false, "23-23=0", 2) | Out-Null
$d.Content.Find.Execute("20+43=63", $true, $false, $false, $false, $false, $true, 1, $false, "49+35=84", 2) | Out-Null
$d.Content.Find.Execute("30+21=51", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=63", 2) | Out-Null
$d.Content.Find.Execute("34+30=64", $true, $false, $false, $false, $false, $true, 1, $false, "41+18=59", 2) | Out-Null
$d.Content.Find.Execute("50-19=31", $true, $false, $false, $false, $false, $true, 1, $false, "13+27=40", 2) | Out-Null
$d.Content.Find.Execute("9+32=41", $true, $false, $false, $false, $false, $true, 1, $false, "57+8=65", 2) | Out-Null
$d.Content.Find.Execute("24+13=37", $true, $false, $false, $false, $false, $true, 1, $false, "15+62=77", 2) | Out-Null
$d.Content.Find.Execute("50-34=16", $true, $false, $false, $false, $false, $true, 1, $false, "19+67=86", 2) | Out-Null
$d.Content.Find.Execute("31+6=37", $true, $false, $false, $false, $false, $true, 1, $false, "96-93=3", 2) | Out-Null
$d.Content.Find.Execute("9+52=61", $true, $false, $false, $false, $false, $true, 1, $false, "40-24=16", 2) | Out-Null
$d.Content.Find.Execute("53-8=45", $true, $false, $false, $false, $false, $true, 1, $false, "5+72=77", 2) | Out-Null
$d.Content.Find.Execute("97-54=43", $true, $false, $false, $false, $false, $true, 1, $false, "77-52=25", 2) | Out-Null
$d.Content.Find.Execute("8-5=3", $true, $false, $false, $false, $false, $true, 1, $false, "6+89=95", 2) | Out-Null
$d.Content.Find.Execute("12+57=69", $true, $false, $false, $false, $false, $true, 1, $false, "96-92=4", 2) | Out-Null
$d.Content.Find.Execute("0+10=10", $true, $false, $false, $false, $false, $true, 1, $false, "77-51=26", 2) | Out-Null
$d.Content.Find.Execute("55+43=98", $true, $false, $false, $false, $false, $true, 1, $false, "83-56=27", 2) | Out-Null
$d.Content.Find.Execute("46+47=93", $true, $false, $false, $false, $false, $true, 1, $false, "20+18=38", 2) | Out-Null
$d.Content.Find.Execute("4+17=21", $true, $false, $false, $false, $false, $true, 1, $false, "40-25=15", 2) | Out-Null
$d.Content.Find.Execute("33-3=30", $true, $false, $false, $false, $false, $true, 1, $false, "6+34=40", 2) | Out-Null
$d.Content.Find.Execute("52+37=89", $true, $false, $false, $false, $false, $true, 1, $false, "25+73=98", 2) | Out-Null
$d.Content.Find.Execute("50+25=75", $true, $false, $false, $false, $false, $true, 1, $false, "94-18=76", 2) | Out-Null
$d.Content.Find.Execute("9+45=54", $true, $false, $false, $false, $false, $true, 1, $false, "26-26=0", 2) | Out-Null
$d.Content.Find.Execute("7+9=16", $true, $false, $false, $false, $false, $true, 1, $false, "89-8=81", 2) | Out-Null
$d.Content.Find.Execute("55-49=6", $true, $false, $false, $false, $false, $true, 1, $false, "48+33=81", 2) | Out-Null
$d.Content.Find.Execute("14+70=84", $true, $false, $false, $false, $false, $true, 1, $false, "4+38=42", 2) | Out-Null
$d.Content.Find.Execute("69-31=38", $true, $false, $false, $false, $false, $true, 1, $false, "50-4=46", 2) | Out-Null
$d.Content.Find.Execute("14+21=35", $true, $false, $false, $false, $false, $true, 1, $false, "84-58=26", 2) | Out-Null
$d.Content.Find.Execute("1+60=61", $true, $false, $false, $false, $false, $true, 1, $false, "64-34=30", 2) | Out-Null
$d.Content.Find.Execute("55-46=9", $true, $false, $false, $false, $false, $true, 1, $false, "75-52=23", 2) | Out-Null
$d.Content.Find.Execute("91-55=36", $true, $false, $false, $false, $false, $true, 1, $false, "38+41=79", 2) | Out-Null
$d.Content.Find.Execute("8+69=77", $true, $false, $false, $false, $false, $true, 1, $false, "1+91=92", 2) | Out-Null
$d.Content.Find.Execute("59-42=17", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=81", 2) | Out-Null
$d.Content.Find.Execute("84-49=35", $true, $false, $false, $false, $false, $true, 1, $false, "40+37=77", 2) | Out-Null
$d.Content.Find.Execute("63-11=52", $true, $false, $false, $false, $false, $true, 1, $false, "0+36=36", 2) | Out-Null
$d.Content.Find.Execute("0+6=6", $true, $false, $false, $false, $false, $true, 1, $false, "10+12=22", 2) | Out-Null
$d.Content.Find.Execute("17+60=77", $true, $false, $false, $false, $false, $true, 1, $false, "55-23=32", 2) | Out-Null
$d.Content.Find.Execute("7+89=96", $true, $false, $false, $false, $false, $true, 1, $false, "9+37=46", 2) | Out-Null
$d.Content.Find.Execute("97-4=93", $true, $false, $false, $false, $false, $true, 1, $false, "10-7=3", 2) | Out-Null
$d.Content.Find.Execute("61+35=96", $true, $false, $false, $false, $false, $true, 1, $false, "89-55=34", 2) | Out-Null
$d.Content.Find.Execute("12+45=57", $true, $false, $false, $false, $false, $true, 1, $false, "87-6=81", 2) | Out-Null
$d.Content.Find.Execute("1+50=51", $true, $false, $false, $false, $false, $true, 1, $false, "23-6=17", 2) | Out-Null
$d.Content.Find.Execute("75+20=95", $true, $false, $false, $false, $false, $true, 1, $false, "15+60=75", 2) | Out-Null
$d.Content.Find.Execute("76-75=1", $true, $false, $false, $false, $false, $true, 1, $false, "73-37=36", 2) | Out-Null
$d.Content.Find.Execute("39+3=42", $true, $false, $false, $false, $false, $true, 1, $false, "23+63=86", 2) | Out-Null
$d.Content.Find.Execute("32-8=24", $true, $false, $false, $false, $false, $true, 1, $false, "81+0=81", 2) | Out-Null
$d.Content.Find.Execute("94+2=96", $true, $false, $false, $false, $false, $true, 1, $false, "55-3=52", 2) | Out-Null
$d.Content.Find.Execute("90-68=22", $true, $false, $false, $false, $false, $true, 1, $false, "6+87=93", 2) | Out-Null
$d.Content.Find.Execute("6+5=11", $true, $false, $false, $false, $false, $true, 1, $false, "31-12=19", 2) | Out-Null
$d.Content.Find.Execute("2+68=70", $true, $false, $false, $false, $false, $true, 1, $false, "51+38=89", 2) | Out-Null
$d.Content.Find.Execute("9+88=97", $true, $false, $false, $false, $false, $true, 1, $false, "64+5=69", 2) | Out-Null
$d.Content.Find.Execute("63-31=32", $true, $false, $false, $false, $false, $true, 1, $false, "84+12=96", 2) | Out-Null
$d.Content.Find.Execute("39-12=27", $true, $false, $false, $false, $false, $true, 1, $false, "87-40=47", 2) | Out-Null
$d.Content.Find.Execute("36+13=49", $true, $false, $false, $false, $false, $true, 1, $false, "19+9=28", 2) | Out-Null
$d.Content.Find.Execute("53-6=47", $true, $false, $false, $false, $false, $true, 1, $false, "27-15=12", 2) | Out-Null
$d.Content.Find.Execute("56+11=67", $true, $false, $false, $false, $false, $true, 1, $false, "92-3=89", 2) | Out-Null
$d.Content.Find.Execute("1+73=74", $true, $false, $false, $false, $false, $true, 1, $false, "69-38=31", 2) | Out-Null
$d.Content.Find.Execute("31-6=25", $true, $false, $false, $false, $false, $true, 1, $false, "91-15=76", 2) | Out-Null
$d.Content.Find.Execute("63-48=15", $true, $false, $false, $false, $false, $true, 1, $false, "80-46=34", 2) | Out-Null
$d.Content.Find.Execute("19+43=62", $true, $false, $false, $false, $false, $true, 1, $false, "30-7=23", 2) | Out-Null
$d.Content.Find.Execute("50-15=35", $true, $false, $false, $false, $false, $true, 1, $false, "81-48=33", 2) | Out-Null
$d.Content.Find.Execute("94-42=52", $true, $false, $false, $false, $false, $true, 1, $false, "13+30=43", 2) | Out-Null
$d.Content.Find.Execute("29+9=38", $true, $false, $false, $false, $false, $true, 1, $false, "72-31=41", 2) | Out-Null
$d.Content.Find.Execute("60+24=84", $true, $false, $false, $false, $false, $true, 1, $false, "51-0=51", 2) | Out-Null
$d.Content.Find.Execute("1+12=13", $true, $false, $false, $false, $false, $true, 1, $false, "5+11=16", 2) | Out-Null
$d.Content.Find.Execute("83-3=80", $true, $false, $false, $false, $false, $true, 1, $false, "5+74=79", 2) | Out-Null
$d.Content.Find.Execute("28-19=9", $true, $false, $false, $false, $false, $true, 1, $false, "93-25=68", 2) | Out-Null
$d.Content.Find.Execute("21+24=45", $true, $false, $false, $false, $false, $true, 1, $false, "86-33=53", 2) | Out-Null
$d.Content.Find.Execute("67+0=67", $true, $false, $false, $false, $false, $true, 1, $false, "74-8=66", 2) | Out-Null
$d.Content.Find.Execute("19-7=12", $true, $false, $false, $false, $false, $true, 1, $false, "92+5=97", 2) | Out-Null
$d.Content.Find.Execute("38-9=29", $true, $false, $false, $false, $false, $true, 1, $false, "75-9=66", 2) | Out-Null
$d.Content.Find.Execute("9+47=56", $true, $false, $false, $false, $false, $true, 1, $false, "70-21=49", 2) | Out-Null
$d.Content.Find.Execute("42+5=47", $true, $false, $false, $false, $false, $true, 1, $false, "55+30=85", 2) | Out-Null
